# Create and check functions to summarize model results
# Rework the lookup-table header labels: split each sheet's "short code" /
# "label" columns so that column A always holds the short "abbr_*" code,
# column B holds the human-readable long-form label (previously "lab_*"),
# and the ordering column is renamed from "ord" to "order".
#
# Sheet "subset"   : A=abbr_subset (was subset), B=subset (was lab_subset), C=order (was ord)
# Sheet "term"     : A=abbr_term   (was term),   B=term   (was lab_term),   C=order (was ord)
# Sheet "mean"     : A=abbr_value  (was value),  B=value  (was lab_value),  C=order (was ord), D=abbr_term (was term)
# Sheet "contrast" : A=abbr_contrast (was contrast), B=contrast (was lab_contrast), C=order (was ord), D=abbr_term (was term)

$wb = $excel.ActiveWorkbook

$wsSubset = $wb.Worksheets.Item("subset")
$wsTerm = $wb.Worksheets.Item("term")
$wsMean = $wb.Worksheets.Item("mean")
$wsContrast = $wb.Worksheets.Item("contrast")

# Set column-A short codes first (introduces the new shared strings in
# "abbr_subset, abbr_term, abbr_value, abbr_contrast" order).
$wsSubset.Range("A1").Value = "abbr_subset"
$wsTerm.Range("A1").Value = "abbr_term"
$wsMean.Range("A1").Value = "abbr_value"
$wsContrast.Range("A1").Value = "abbr_contrast"

# Column-B long labels now reuse the already-present short-name strings.
$wsSubset.Range("B1").Value = "subset"
$wsTerm.Range("B1").Value = "term"
$wsMean.Range("B1").Value = "value"
$wsContrast.Range("B1").Value = "contrast"

# Rename the ordering column ("ord" -> "order") last, introducing "order" as
# the final new shared string.
$wsSubset.Range("C1").Value = "order"
$wsTerm.Range("C1").Value = "order"
$wsMean.Range("C1").Value = "order"
$wsContrast.Range("C1").Value = "order"

# The "mean" and "contrast" sheets also carry a term lookup column that is
# renamed the same way as the "term" sheet's own header.
$wsMean.Range("D1").Value = "abbr_term"
$wsContrast.Range("D1").Value = "abbr_term"

# Reset the selection on each sheet to C1, matching the post-edit view state.
$wsSubset.Range("C1").Select() | Out-Null
$wsTerm.Range("C1").Select() | Out-Null
$wsMean.Range("C1").Select() | Out-Null
$wsContrast.Range("C1").Select() | Out-Null

$wsSubset.Activate() | Out-Null
